# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns,
# and for a few rows (39-40, 48-51) the Coin (B) and Link (C) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.402.50"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").Value = "1.847.61"
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").Value = "'240.49"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("D6").Value = "'0.6273"
$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("D7").Value = "'0.9990"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("E8").Value = "  -1.65%  "

$ws.Range("D9").Value = "'0.2907"
$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").Value = "'24.46"
$ws.Range("E10").Value = "  -1.12%  "

$ws.Range("D11").Value = "'0.07738"
$ws.Range("E11").Value = "  -0.05%  "

$ws.Range("D12").Value = "1.847.27"
$ws.Range("E12").Value = "  -2.17%  "

$ws.Range("D13").Value = "'5.002"
$ws.Range("E13").Value = "  -0.67%  "

$ws.Range("D14").Value = "'0.6811"
$ws.Range("E14").Value = "  +0.35%  "

$ws.Range("D15").Value = "'0.00001057"
$ws.Range("E15").Value = "  +0.35%  "

$ws.Range("D16").Value = "'82.23"
$ws.Range("E16").Value = "  -1.20%  "

$ws.Range("D17").Value = "2.108.28"
$ws.Range("E17").Value = "  -3.60%  "

$ws.Range("D18").Value = "'6.181"
$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("D19").Value = "29.429.39"
$ws.Range("E19").Value = "  +0.06%  "

$ws.Range("D20").Value = "'229.77"
$ws.Range("E20").Value = "  +0.86%  "

$ws.Range("D21").Value = "'12.34"
$ws.Range("E21").Value = "  -0.19%  "

$ws.Range("D22").Value = "'0.9992"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").Value = "'7.490"
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").Value = "'0.9997"
$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("D25").Value = "'159.13"
$ws.Range("E25").Value = "  +0.27%  "

$ws.Range("D26").Value = "'0.1377"
$ws.Range("E26").Value = "  -0.77%  "

$ws.Range("D27").Value = "'8.413"
$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D28").Value = "'17.54"
$ws.Range("E28").Value = "  -0.84%  "

$ws.Range("D29").Value = "'0.06468"
$ws.Range("E29").Value = "  +15.43%  "

$ws.Range("D30").Value = "'1.418"
$ws.Range("E30").Value = "  +2.14%  "

$ws.Range("D31").Value = "'1.476"
$ws.Range("E31").Value = "  +1.16%  "

$ws.Range("D32").Value = "'4.099"
$ws.Range("E32").Value = "  -0.27%  "

$ws.Range("D33").Value = "'4.100"
$ws.Range("E33").Value = "  +0.78%  "

$ws.Range("D34").Value = "'1.832"
$ws.Range("E34").Value = "  -0.33%  "

$ws.Range("E35").Value = "  -1.80%  "

$ws.Range("D36").Value = "'0.6968"
$ws.Range("E36").Value = "  -0.62%  "

$ws.Range("D37").Value = "'2.579"
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("D38").Value = "1.272.46"
$ws.Range("E38").Value = "  +3.22%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.824"
$ws.Range("E39").Value = "  +4.04%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01834"
$ws.Range("E40").Value = "  +1.50%  "

$ws.Range("D41").Value = "'6.758"
$ws.Range("E41").Value = "  +5.80%  "

$ws.Range("D42").Value = "'0.9087"
$ws.Range("E42").Value = "  +0.68%  "

$ws.Range("D43").Value = "'0.9990"

$ws.Range("D44").Value = "2.011.00"
$ws.Range("E44").Value = "  -18.28%  "

$ws.Range("D45").Value = "'101.44"
$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("D46").Value = "'66.42"
$ws.Range("E46").Value = "  +0.46%  "

$ws.Range("D47").Value = "'1.747"
$ws.Range("E47").Value = "  +3.98%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.00000000119"
$ws.Range("E48").Value = "  -0.41%  "

$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'7.085"
$ws.Range("E49").Value = "  -1.83%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.1176"
$ws.Range("E50").Value = "  +3.86%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'9.074"
$ws.Range("E51").Value = "  +0.50%  "
